$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.026972333333333
$ws.Range("H2").Value = 3.080917
$ws.Range("I2").Value = 0.2032541865322035
$ws.Range("J2").Value = 0.2032541865322035
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.316186
$ws.Range("N2").Value = 9.948558
$ws.Range("O2").Value = 0.7638813129544791
$ws.Range("P2").Value = 0.7638813129544793
$ws.Range("Q2").Value = 3.405631274187333
$ws.Range("R2").Value = 30.650681467686
$ws.Range("S2").Value = 0.1552620748717142
$ws.Range("T2").Value = 0.1552620748717142
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.026972333333333
$ws.Range("H3").Value = 3.080917
$ws.Range("I3").Value = 0.2032541865322035
$ws.Range("J3").Value = 0.2032541865322035
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.4282866666666667
$ws.Range("N3").Value = 1.28486
$ws.Range("O3").Value = 0.09865555829927233
$ws.Range("P3").Value = 0.09865555829927235
$ws.Range("Q3").Value = 0.4398385574022223
$ws.Range("R3").Value = 3.95854701662
$ws.Range("S3").Value = 0.02005215524899897
$ws.Range("T3").Value = 0.02005215524899898
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.026972333333333
$ws.Range("H4").Value = 3.080917
$ws.Range("I4").Value = 0.2032541865322035
$ws.Range("J4").Value = 0.2032541865322035
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.417289
$ws.Range("N4").Value = 1.251867
$ws.Range("O4").Value = 0.09612225285356782
$ws.Range("P4").Value = 0.09612225285356783
$ws.Range("Q4").Value = 0.4285442580043334
$ws.Range("R4").Value = 3.856898322039
$ws.Range("S4").Value = 0.0195372503113947
$ws.Range("T4").Value = 0.0195372503113947
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.026972333333333
$ws.Range("H5").Value = 3.080917
$ws.Range("I5").Value = 0.2032541865322035
$ws.Range("J5").Value = 0.2032541865322035
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.1794703333333333
$ws.Range("N5").Value = 0.538411
$ws.Range("O5").Value = 0.04134087589268053
$ws.Range("P5").Value = 0.04134087589268053
$ws.Range("Q5").Value = 0.1843110669874444
$ws.Range("R5").Value = 1.658799602887
$ws.Range("S5").Value = 0.008402706100095563
$ws.Range("T5").Value = 0.008402706100095563
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.358031333333334
$ws.Range("H6").Value = 4.074094000000001
$ws.Range("I6").Value = 0.2687760370778347
$ws.Range("J6").Value = 0.2687760370778347
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.316186
$ws.Range("N6").Value = 9.948558
$ws.Range("O6").Value = 0.7638813129544791
$ws.Range("P6").Value = 0.7638813129544793
$ws.Range("Q6").Value = 4.503484495161334
$ws.Range("R6").Value = 40.53136045645201
$ws.Range("S6").Value = 0.2053129920937181
$ws.Range("T6").Value = 0.2053129920937181
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.358031333333334
$ws.Range("H7").Value = 4.074094000000001
$ws.Range("I7").Value = 0.2687760370778347
$ws.Range("J7").Value = 0.2687760370778347
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.4282866666666667
$ws.Range("N7").Value = 1.28486
$ws.Range("O7").Value = 0.09865555829927233
$ws.Range("P7").Value = 0.09865555829927235
$ws.Range("Q7").Value = 0.5816267129822223
$ws.Range("R7").Value = 5.234640416840001
$ws.Range("S7").Value = 0.0265162499953797
$ws.Range("T7").Value = 0.0265162499953797
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.358031333333334
$ws.Range("H8").Value = 4.074094000000001
$ws.Range("I8").Value = 0.2687760370778347
$ws.Range("J8").Value = 0.2687760370778347
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.417289
$ws.Range("N8").Value = 1.251867
$ws.Range("O8").Value = 0.09612225285356782
$ws.Range("P8").Value = 0.09612225285356783
$ws.Range("Q8").Value = 0.5666915370553335
$ws.Range("R8").Value = 5.100223833498001
$ws.Range("S8").Value = 0.02583535819697554
$ws.Range("T8").Value = 0.02583535819697555
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.358031333333334
$ws.Range("H9").Value = 4.074094000000001
$ws.Range("I9").Value = 0.2687760370778347
$ws.Range("J9").Value = 0.2687760370778347
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.1794703333333333
$ws.Range("N9").Value = 0.538411
$ws.Range("O9").Value = 0.04134087589268053
$ws.Range("P9").Value = 0.04134087589268053
$ws.Range("Q9").Value = 0.2437263360704445
$ws.Range("R9").Value = 2.193537024634
$ws.Range("S9").Value = 0.01111143679176126
$ws.Range("T9").Value = 0.01111143679176126
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.667646666666666
$ws.Range("H10").Value = 8.002939999999999
$ws.Range("I10").Value = 0.5279697763899619
$ws.Range("J10").Value = 0.5279697763899619
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.316186
$ws.Range("N10").Value = 9.948558
$ws.Range("O10").Value = 0.7638813129544791
$ws.Range("P10").Value = 0.7638813129544793
$ws.Range("Q10").Value = 8.846412528946665
$ws.Range("R10").Value = 79.61771276051999
$ws.Range("S10").Value = 0.4033062459890469
$ws.Range("T10").Value = 0.4033062459890469
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.667646666666666
$ws.Range("H11").Value = 8.002939999999999
$ws.Range("I11").Value = 0.5279697763899619
$ws.Range("J11").Value = 0.5279697763899619
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.4282866666666667
$ws.Range("N11").Value = 1.28486
$ws.Range("O11").Value = 0.09865555829927233
$ws.Range("P11").Value = 0.09865555829927235
$ws.Range("Q11").Value = 1.142517498711111
$ws.Range("R11").Value = 10.2826574884
$ws.Range("S11").Value = 0.05208715305489366
$ws.Range("T11").Value = 0.05208715305489367
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.667646666666666
$ws.Range("H12").Value = 8.002939999999999
$ws.Range("I12").Value = 0.5279697763899619
$ws.Range("J12").Value = 0.5279697763899619
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.417289
$ws.Range("N12").Value = 1.251867
$ws.Range("O12").Value = 0.09612225285356782
$ws.Range("P12").Value = 0.09612225285356783
$ws.Range("Q12").Value = 1.113179609886667
$ws.Range("R12").Value = 10.01861648898
$ws.Range("S12").Value = 0.05074964434519758
$ws.Range("T12").Value = 0.05074964434519758
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.667646666666666
$ws.Range("H13").Value = 8.002939999999999
$ws.Range("I13").Value = 0.5279697763899619
$ws.Range("J13").Value = 0.5279697763899619
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.1794703333333333
$ws.Range("N13").Value = 0.538411
$ws.Range("O13").Value = 0.04134087589268053
$ws.Range("P13").Value = 0.04134087589268053
$ws.Range("Q13").Value = 0.4787634364822221
$ws.Range("R13").Value = 4.308870928339999
$ws.Range("S13").Value = 0.02182673300082371
$ws.Range("T13").Value = 0.02182673300082371
